# Append the new Argent price row (2025-01-21, 42.6) to Sheet1, matching
# the existing "date"/"Value" columns which are stored as plain text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 84

# Force text formatting before assigning, otherwise Excel auto-detects
# "2025-01-21" as a date and "42.6" could be read back as a number.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2025-01-21"

$ws.Cells.Item($newRow, 2).NumberFormat = "@"
$ws.Cells.Item($newRow, 2).Value = "42.6"
